$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 2 first (so the first part's character offsets, which come
# earlier in the document, are not disturbed by this edit):
# Merge "Estimate " + (bookmark _GoBack) + "CPMG phase map: " into a
# single run "Estimate CPMG phase map: " and drop the old bookmark.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Estimate CPMG phase map: ")
$start2 = $r2.Start
$end2 = $r2.End

$mergeRange = $d.Range($start2, $end2)
# Force an actual text change (so the engine really rewrites the run,
# merging the two original runs and dropping the bookmark that sat
# between them), then fix the text back to the desired final value.
$mergeRange.Text = "Estimate CPMG phase map: #TMP#"

$fixRange = $d.Content
$fixRange.Find.Execute("Estimate CPMG phase map: #TMP#")
$fixRange.Text = "Estimate CPMG phase map: "

# ---------------------------------------------------------------------
# Part 1: split "./data: the folder contains test data for recon: ...
# ... where " into "./data: the folder contains test data" + new
# bookmark _GoBack + " (R2) for recon: ... where ".
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("./data: the folder contains test data for recon: haste (DW-HASTE) and sp (SP-DW-HASTE). where ")
$start1 = $r1.Start
$end1 = $r1.End

# Replace the whole run with two anchor placeholder characters so the
# run's original rPr / run boundaries with neighbouring runs are kept
# intact (a fully-emptied range loses its run and subsequent inserts
# pick up the wrong/neighbouring formatting).
$full1 = $d.Range($start1, $end1)
$full1.Text = "#A##B#"

# Replace the first placeholder with the first piece of text.
$partA = $d.Range($start1, $start1 + 3)
$partA.Text = "./data: the folder contains test data"

$splitPos = $start1 + ("./data: the folder contains test data").Length

# Insert the (new) _GoBack bookmark exactly at the split point, before
# touching the second placeholder, so it lands between the two runs.
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Replace the second placeholder with the remaining text. Because this
# is a genuine text replacement (not an insert at a point), the
# bookmark added above stays put instead of being dragged along.
$partB = $d.Range($splitPos, $splitPos + 3)
$partB.Text = " (R2) for recon: haste (DW-HASTE) and sp (SP-DW-HASTE). where "
